$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.919.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.254.54'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.648'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.29%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '231.15'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.23'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.99%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.452'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0981'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.06'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.78'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +11.79%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.587.06'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.51'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.10'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.830'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.253.52'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.840.69'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.20'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.05'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.55'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.65%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.35'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +23.63%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.83'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.95'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.126'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0703'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +6.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.81'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.88'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.59%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.46%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.16%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.82%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.000220'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.43'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.49%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.23'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.69%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.47'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '97.51'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.03%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.19'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0944'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.36'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.437.54'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.95'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.56%  '
